$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- Sheet "Elements": swap the AK and AL columns (header + all data rows) ---
$elements = $wb.Worksheets.Item("Elements")

$lastRow = 16
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value()
    $alVal = $alCell.Value()
    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the column widths to match (AK was 24.98046875, AL was 80.21875 -> now swapped)
# (ColumnWidth setter on this host quantizes to 1/6-character steps, so these are
# the closest achievable inputs to the exact target widths of 80.21875 / 24.98046875)
$elements.Columns.Item(37).ColumnWidth = 79.33333333333333
$elements.Columns.Item(38).ColumnWidth = 24.166666666666668
